$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the extraneous leading index column (old column A). This shifts the
# real data (old B:F) left into A:E, carrying values AND formatting with it
# (row 1 keeps its bold/bordered header style, data rows lose the stray
# bold/bordered style that had wrongly been applied to column A).
$ws.Columns.Item(1).Delete()

# Fix header typo: MODEL_CONDITION -> MODELCONDITION (now in column D).
$ws.Range("D1").Value = "MODELCONDITION"
